# feat: window halaman utama
# Adds a new row (row 3) to the active sheet with values "a" / "1",
# extending the used range from A1:B2 to A1:B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "a"

# "1" must stay a text value (not be coerced to the number 1), so force it
# with a leading apostrophe like a user typing into the cell would, then
# reset the style back to Normal so no stray quote-prefix style lingers.
$ws.Range("B3").Value = "'1"
$ws.Range("B3").Style = "Normal"
